# Apply the "GitHub Actions" cryptos-list refresh (Sun Jan  7 18:27:26 UTC 2024).
# Every data cell in the sheet is stored as text (inlineStr), including the
# numeric-looking Price column, so plain numeric-looking values are written with
# a leading apostrophe to force Excel's text interpretation (matches how a user
# typing a number-like value into a pre-formatted text cell behaves), keeping the
# stored cell content an exact string match without the apostrophe itself.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "44.581.86"
$ws.Range("E2").Value = "  +0.62%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "2.242.02"
$ws.Range("E3").Value = "  -0.35%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.32%  "

# Row 5 (BNB)
$ws.Range("D5").Value = "'305.53"
$ws.Range("E5").Value = "  -0.73%  "

# Row 6 (Solana)
$ws.Range("D6").Value = "'94.97"
$ws.Range("E6").Value = "  -1.82%  "

# Row 7 (XRP)
$ws.Range("D7").Value = "'0.571"
$ws.Range("E7").Value = "  -0.70%  "

# Row 8 (USDC)
$ws.Range("E8").Value = "  +0.18%  "

# Row 9 (Cardano)
$ws.Range("D9").Value = "'0.519"
$ws.Range("E9").Value = "  -1.92%  "

# Row 10 (Avalanche)
$ws.Range("D10").Value = "'35.44"
$ws.Range("E10").Value = "  +0.47%  "

# Row 11 (Dogecoin)
$ws.Range("D11").Value = "'0.0803"
$ws.Range("E11").Value = "  -1.72%  "

# Row 12 (Polkadot)
$ws.Range("D12").Value = "'7.21"
$ws.Range("E12").Value = "  -1.26%  "

# Row 13 (TRON)
$ws.Range("E13").Value = "  -0.38%  "

# Row 14 (WrappedliquidstakedEther2.0)
$ws.Range("D14").Value = "2.584.10"

# Row 15 (WrappedEther)
$ws.Range("D15").Value = "2.237.17"
$ws.Range("E15").Value = "  +0.05%  "

# Row 16 (Polygon)
$ws.Range("D16").Value = "'0.833"
$ws.Range("E16").Value = "  -0.48%  "

# Row 17 (Chainlink)
$ws.Range("D17").Value = "'13.56"
$ws.Range("E17").Value = "  -0.71%  "

# Row 18 (WrappedBTC)
$ws.Range("D18").Value = "44.397.52"
$ws.Range("E18").Value = "  +0.67%  "

# Row 19 (ShibaInu)
$ws.Range("D19").Value = "0.0₃0940"
$ws.Range("E19").Value = "  -3.27%  "

# Row 20 (InternetComputer(DFINITY))
$ws.Range("D20").Value = "'11.85"
$ws.Range("E20").Value = "  -3.07%  "

# Row 21 (Uniswap)
$ws.Range("E21").Value = "  -3.43%  "

# Row 22 (Litecoin)
$ws.Range("D22").Value = "'65.28"
$ws.Range("E22").Value = "  -0.78%  "

# Row 23 (BitcoinCash)
$ws.Range("D23").Value = "'237.22"
$ws.Range("E23").Value = "  -0.34%  "

# Row 24 (PancakeSwap)
$ws.Range("D24").Value = "'2.94"
$ws.Range("E24").Value = "  -0.88%  "

# Row 25 (ImmutableX)
$ws.Range("D25").Value = "'1.99"
$ws.Range("E25").Value = "  -1.92%  "

# Row 26 (Dai)
$ws.Range("E26").Value = "  -0.17%  "

# Row 27 (Toncoin)
$ws.Range("E27").Value = "  +8.02%  "

# Row 28 (Cosmos)
$ws.Range("D28").Value = "'9.75"
$ws.Range("E28").Value = "  -2.51%  "

# Row 29 (InjectiveProtocol)
$ws.Range("D29").Value = "'37.10"
$ws.Range("E29").Value = "  -4.50%  "

# Row 30 (Filecoin)
$ws.Range("D30").Value = "'5.91"
$ws.Range("E30").Value = "  -0.72%  "

# Row 31 (EthereumClassic)
$ws.Range("D31").Value = "'19.89"
$ws.Range("E31").Value = "  -1.29%  "

# Row 32 (Monero)
$ws.Range("D32").Value = "'150.03"
$ws.Range("E32").Value = "  -1.71%  "

# Row 33 (Hedera)
$ws.Range("D33").Value = "'0.0785"
$ws.Range("E33").Value = "  -1.84%  "

# Row 34 (WEMIXToken)
$ws.Range("E34").Value = "  -0.02%  "

# Row 35 (LidoDAOToken)
$ws.Range("D35").Value = "'3.15"
$ws.Range("E35").Value = "  -3.19%  "

# Row 36 (Kaspa)
$ws.Range("E36").Value = "  +0.35%  "

# Row 37 (ARBITRUM)
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").Value = "'0.118"
$ws.Range("E37").Value = "  -1.76%  "

# Row 38 (Stellar)
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'1.87"
$ws.Range("E38").Value = "  +5.93%  "

# Row 39 (Celestia)
$ws.Range("D39").Value = "'15.15"
$ws.Range("E39").Value = "  +2.71%  "

# Row 40 (NEARProtocol)
$ws.Range("D40").Value = "'3.39"
$ws.Range("E40").Value = "  -7.56%  "

# Row 41 (RenderToken)
$ws.Range("D41").Value = "'3.79"
$ws.Range("E41").Value = "  -1.99%  "

# Row 42 (VeChain)
$ws.Range("D42").Value = "'0.0300"
$ws.Range("E42").Value = "  -0.30%  "

# Row 43 (FirstDigitalUSD)
$ws.Range("E43").Value = "  +0.15%  "

# Row 44 (Maker)
$ws.Range("D44").Value = "1.810.49"
$ws.Range("E44").Value = "  +2.94%  "

# Row 45 (Stacks)
$ws.Range("D45").Value = "'1.78"
$ws.Range("E45").Value = "  +11.56%  "

# Row 46 (BitcoinSV)
$ws.Range("D46").Value = "'81.12"
$ws.Range("E46").Value = "  -2.66%  "

# Row 47 (Algorand)
$ws.Range("E47").Value = "  -1.94%  "

# Row 48 (Aave)
$ws.Range("D48").Value = "'98.30"
$ws.Range("E48").Value = "  -2.41%  "

# Row 49 (THORChain)
$ws.Range("D49").Value = "'4.84"
$ws.Range("E49").Value = "  -2.80%  "

# Row 50 (ordi)
$ws.Range("D50").Value = "'68.43"
$ws.Range("E50").Value = "  +0.46%  "

# Row 51 (MultiversX)
$ws.Range("D51").Value = "'54.13"
$ws.Range("E51").Value = "  -1.59%  "
